$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 (I0) and J1 (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match formatting (bold font, thin border, centered/top alignment) of the existing header row
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Data rows: row, I0 value, IF value
$data = @(
    @(2,4,4),
    @(3,5,5),
    @(4,6,7),
    @(5,6,6),
    @(6,9,9),
    @(7,8,8),
    @(8,7,7),
    @(9,4,4),
    @(10,6,6),
    @(11,6,6),
    @(12,9,9),
    @(13,6,7),
    @(14,7,7),
    @(15,6,6),
    @(16,6,6),
    @(17,8,8),
    @(18,8,8),
    @(19,5,6),
    @(20,9,9),
    @(21,5,6),
    @(22,6,7),
    @(23,4,5),
    @(24,9,9),
    @(25,10,10),
    @(26,9,9),
    @(27,8,8),
    @(28,6,6),
    @(29,8,8),
    @(30,9,9),
    @(31,8,8),
    @(32,8,8),
    @(33,8,8),
    @(34,7,7),
    @(35,7,7),
    @(36,7,7),
    @(37,3,3),
    @(38,7,7),
    @(39,9,9),
    @(40,7,8),
    @(41,7,8),
    @(42,8,8),
    @(43,8,9),
    @(44,7,8),
    @(45,8,8),
    @(46,9,9),
    @(47,8,8),
    @(48,9,9),
    @(49,7,7),
    @(50,6,6),
    @(51,6,8),
    @(52,8,8),
    @(53,5,5),
    @(54,7,7),
    @(55,10,10),
    @(56,10,10),
    @(57,6,7),
    @(58,5,7),
    @(59,7,7),
    @(60,7,7),
    @(61,5,5),
    @(62,8,8),
    @(63,7,7),
    @(64,7,7),
    @(65,6,6),
    @(66,5,5),
    @(67,3,3)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
